$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 131136941
$ws.Range("B4").Value = 83090
$ws.Range("E4").Value = 1312
$ws.Range("F4").Value = "Gammelgransskål"
$ws.Range("G4").Value = "Pseudographis pinicola"
$ws.Range("H4").Value = "(Nyl.) Rehm"
$ws.Range("J4").Value = "fruktkroppar"
$ws.Range("Q4").Value = 788995
$ws.Range("R4").Value = 7131220
$ws.Range("AC4").Value = "på en gammal senvuxen gran"

# Row 5
$ws.Range("A5").Value = 131136874
$ws.Range("B5").Value = 79244
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("J5").Value = "bålar"
$ws.Range("Q5").Value = 788960
$ws.Range("R5").Value = 7131416
$ws.Range("AC5").Value = ""
$ws.Range("AE5").Value = $false

# Row 6
$ws.Range("A6").Value = 131136961
$ws.Range("B6").Value = 57884
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("J6").Value = ""
$ws.Range("Q6").Value = 789068
$ws.Range("R6").Value = 7131245
$ws.Range("AC6").Value = "barksprätt på gammal gran"
$ws.Range("AE6").Value = $true

# Row 7
$ws.Range("A7").Value = 131136881
$ws.Range("B7").Value = 79244
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = "Garnlav"
$ws.Range("G7").Value = "Alectoria sarmentosa"
$ws.Range("H7").Value = "(Ach.) Ach."
$ws.Range("J7").Value = "bålar"
$ws.Range("M7").Value = ""
$ws.Range("Q7").Value = 788972
$ws.Range("R7").Value = 7131396
$ws.Range("AC7").Value = ""

# Row 8
$ws.Range("A8").Value = 131136984
$ws.Range("B8").Value = 57881
$ws.Range("E8").Value = 100049
$ws.Range("F8").Value = "Spillkråka"
$ws.Range("G8").Value = "Dryocopus martius"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("J8").Value = ""
$ws.Range("M8").Value = "äldre spår"
$ws.Range("Q8").Value = 788839
$ws.Range("R8").Value = 7131504
$ws.Range("AC8").Value = "hål i tallstam"
